$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GanttChart")
$ws.Range("H4").Value = 10
